$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 0.0006944444444444445
$ws.Range("K2").Value = 4974
$ws.Range("L2").Value = 0.009948
